$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6720.6665
$ws.Range("I62").Value = 6720.6665
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 6720.6665
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -6096.6665
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 6720.6665
$ws.Range("I65").Value = 6720.6665
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 33603.3325
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -30483.3325
$ws.Range("N65").ClearContents()
$ws.Range("H94").Value = 1258.4615
$ws.Range("I94").Value = 1221.6666
$ws.Range("K94").Value = 1221.6666
$ws.Range("M94").Value = -770.6666
$ws.Range("H107").Value = 1648.3513
$ws.Range("I107").Value = 1376.2593
$ws.Range("J107").Value = 2383
$ws.Range("K107").Value = 1376.2593
$ws.Range("L107").Value = 2383
$ws.Range("M107").Value = 543.7407000000001
$ws.Range("N107").Value = -6223
$ws.Range("H129").Value = 2937.276
$ws.Range("I129").Value = 1550
$ws.Range("K129").Value = 4650
$ws.Range("M129").Value = 350
$ws.Range("H132").Value = 4968.3516
$ws.Range("I132").Value = 4277.4
$ws.Range("K132").Value = 12832.2
$ws.Range("M132").Value = -10302.2
$ws.Range("H135").Value = 1015.7931
$ws.Range("I135").Value = 923.38464
$ws.Range("K135").Value = 8310.46176
$ws.Range("M135").Value = -5775.46176
$ws.Range("H137").Value = 3368.6667
$ws.Range("I137").Value = 3491
$ws.Range("J137").Value = 3058.4644
$ws.Range("K137").Value = 10473
$ws.Range("L137").Value = 9175.393199999999
$ws.Range("M137").Value = -7923
$ws.Range("N137").Value = -14275.3932
$ws.Range("H138").Value = 2897.6042
$ws.Range("I138").Value = 2191
$ws.Range("J138").Value = 3402.3215
$ws.Range("K138").Value = 6573
$ws.Range("L138").Value = 10206.9645
$ws.Range("M138").Value = -1433
$ws.Range("N138").Value = -20486.9645
$ws.Range("H141").Value = 1385.9375
$ws.Range("I141").Value = 1466.9231
$ws.Range("K141").Value = 4400.7693
$ws.Range("M141").Value = 779.2307000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4465.5
$ws.Range("I61").Value = 2334.9583
$ws.Range("K61").Value = 2334.9583
$ws.Range("M61").Value = -2122.9583
$ws.Range("H74").Value = 3076.4443
$ws.Range("I74").Value = 1178.2
$ws.Range("J74").Value = 7390.636
$ws.Range("K74").Value = 1178.2
$ws.Range("L74").Value = 7390.636
$ws.Range("M74").Value = -304.2
$ws.Range("N74").Value = -9138.636
$ws.Range("H77").Value = 3076.4443
$ws.Range("I77").Value = 1178.2
$ws.Range("J77").Value = 7390.636
$ws.Range("K77").Value = 5891
$ws.Range("L77").Value = 36953.18
$ws.Range("M77").Value = -1523
$ws.Range("N77").Value = -45689.18
$ws.Range("H110").Value = 2701.1667
$ws.Range("I110").Value = 2491.7
$ws.Range("J110").Value = 3748.5
$ws.Range("K110").Value = 2491.7
$ws.Range("L110").Value = 3748.5
$ws.Range("M110").Value = -446.6999999999998
$ws.Range("N110").Value = -7838.5
$ws.Range("H136").Value = 4465.5
$ws.Range("I136").Value = 2334.9583
$ws.Range("K136").Value = 7004.874899999999
$ws.Range("M136").Value = -4454.874899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 714.8333
$ws.Range("I94").Value = 684.41174
$ws.Range("J94").Value = 788.7143
$ws.Range("K94").Value = 684.41174
$ws.Range("L94").Value = 788.7143
$ws.Range("M94").Value = -233.41174
$ws.Range("N94").Value = -1690.7143
$ws.Range("H107").Value = 999.6667
$ws.Range("I107").Value = 894.4737
$ws.Range("K107").Value = 894.4737
$ws.Range("M107").Value = 1025.5263

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1571.8334
$ws.Range("I31").Value = 1335.9
$ws.Range("J31").Value = 1662.5769
$ws.Range("K31").Value = 1335.9
$ws.Range("L31").Value = 1662.5769
$ws.Range("M31").Value = -1040.9
$ws.Range("N31").Value = -2252.5769
$ws.Range("H34").Value = 1571.8334
$ws.Range("I34").Value = 1335.9
$ws.Range("J34").Value = 1662.5769
$ws.Range("K34").Value = 1335.9
$ws.Range("L34").Value = 1662.5769
$ws.Range("M34").Value = -1133.9
$ws.Range("N34").Value = -2066.5769
$ws.Range("H134").Value = 8081.75
$ws.Range("I134").Value = 7754.6523
$ws.Range("K134").Value = 23263.9569
$ws.Range("M134").Value = -20728.9569

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2463.182
$ws.Range("J68").Value = 2566.2222
$ws.Range("L68").Value = 7698.6666
$ws.Range("N68").Value = -9320.6666
$ws.Range("H71").Value = 2463.182
$ws.Range("J71").Value = 2566.2222
$ws.Range("L71").Value = 23095.9998
$ws.Range("N71").Value = -31207.9998
$ws.Range("H140").Value = 1939.5758
$ws.Range("I140").Value = 1845.1724
$ws.Range("K140").Value = 5535.5172
$ws.Range("M140").Value = -355.5172000000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 16604.562
$ws.Range("I13").Value = 16902.166
$ws.Range("J13").Value = 15711.75
$ws.Range("K13").Value = 16902.166
$ws.Range("L13").Value = 15711.75
$ws.Range("M13").Value = -16763.166
$ws.Range("N13").Value = -15989.75
$ws.Range("H122").Value = 2084.9167
$ws.Range("I122").Value = 1398.6316
$ws.Range("J122").Value = 4692.8
$ws.Range("K122").Value = 4195.8948
$ws.Range("L122").Value = 14078.4
$ws.Range("M122").Value = -1745.8948
$ws.Range("N122").Value = -18978.4
$ws.Range("H126").Value = 1749.1666
$ws.Range("I126").Value = 1799
$ws.Range("J126").Value = 1649.5
$ws.Range("K126").Value = 5397
$ws.Range("L126").Value = 4948.5
$ws.Range("M126").Value = -2927
$ws.Range("N126").Value = -9888.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3774
$ws.Range("I122").Value = 3615
$ws.Range("J122").Value = 5205
$ws.Range("K122").Value = 10845
$ws.Range("L122").Value = 15615
$ws.Range("M122").Value = -8395
$ws.Range("N122").Value = -20515

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 7500
$ws.Range("J45").Value = 7500
$ws.Range("L45").Value = 7500
$ws.Range("N45").Value = -8482
$ws.Range("H132").Value = 14513.712
$ws.Range("I132").Value = 10992.574
$ws.Range("K132").Value = 32977.722
$ws.Range("M132").Value = -30447.722
